$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new data row at row 28 ("جل صبار برطمان") - pushes the
# existing rows 28 ("سرنجات 3 سم"), 29 ("سرنجات 5 سم"), 30 ("شاش 5سم"),
# 31 (totals) and 32 (footer) down to 29, 30, 31, 32 and 33.
# ------------------------------------------------------------------
$ws.Rows("28:28").Insert()

# Merge the cells of the new row the same way every other item row is
# merged (A:B, C:G, H:K, L:M, N:O - P and Q stay single cells).
$ws.Range("A28:B28").Merge()
$ws.Range("C28:G28").Merge()
$ws.Range("H28:K28").Merge()
$ws.Range("L28:M28").Merge()
$ws.Range("N28:O28").Merge()

# Write the values for the new item. Force text number-format first so
# the numeric-looking strings are stored as text (matching every other
# row in this table) instead of being coerced into real numbers.
$ws.Range("A28").Value = 22

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "جل صبار برطمان"

$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "1:0"

$ws.Range("L28").NumberFormat = "@"
$ws.Range("L28").Value = "0"

$ws.Range("N28").NumberFormat = "@"
$ws.Range("N28").Value = "25.00"

$ws.Range("P28").NumberFormat = "@"
$ws.Range("P28").Value = "25.0000"

$ws.Range("Q28").NumberFormat = "@"
$ws.Range("Q28").Value = "1:0"

# Re-apply the formatting (fonts/fills/borders/number-formats) of the
# row below (the item that used to be row 28) onto the new row so the
# new row looks identical to every other item row.
$ws.Range("A29:Q29").Copy()
$ws.Range("A28:Q28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The row-height isn't carried over by Insert()/PasteSpecial, restore it.
$ws.Rows("28:28").RowHeight = 24.75

# ------------------------------------------------------------------
# Update the grand-total cell (now row 32, used to be row 31): add the
# price of the new item (25.00) to the previous total (1336.17).
# ------------------------------------------------------------------
$ws.Range("P32").Value = 1361.17

# ------------------------------------------------------------------
# Update the footer timestamp (now row 33, used to be row 32).
# ------------------------------------------------------------------
$ws.Range("A33").Value = "Monday, 29 September, 2025 11:42 AM"
